$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update YearlyEnergy (column K) values for rows 3-62
$ws.Range("K3").Value = 333510.13
$ws.Range("K4").Value = 144010.17000000001
$ws.Range("K5").Value = 568353.80000000005
$ws.Range("K6").Value = 31969.133000000002
$ws.Range("K7").Value = 95746.48
$ws.Range("K8").Value = 16722.208999999999
$ws.Range("K9").Value = 20176.07
$ws.Range("K10").Value = 27227.153999999999
$ws.Range("K11").Value = 206903.03
$ws.Range("K12").Value = 377165
$ws.Range("K13").Value = 217364.14
$ws.Range("K14").Value = 2501681.5
$ws.Range("K15").Value = 38155.207000000002
$ws.Range("K16").Value = 15053.463
$ws.Range("K17").Value = 251645.22
$ws.Range("K18").Value = 31206.947
$ws.Range("K19").Value = 178301.69
$ws.Range("K20").Value = 27977.474999999999
$ws.Range("K21").Value = 275806.13
$ws.Range("K22").Value = 676217.25
$ws.Range("K23").Value = 75689.97
$ws.Range("K24").Value = 151941.29999999999
$ws.Range("K25").Value = 14786.965
$ws.Range("K26").Value = 17722.396000000001
$ws.Range("K27").Value = 28764.782999999999
$ws.Range("K28").Value = 120153.89
$ws.Range("K29").Value = 26889.47
$ws.Range("K30").Value = 43730.133000000002
$ws.Range("K31").Value = 64794.336000000003
$ws.Range("K32").Value = 504351
$ws.Range("K33").Value = 63327.042999999998
$ws.Range("K34").Value = 122112.23
$ws.Range("K35").Value = 35462.75
$ws.Range("K36").Value = 82257.33
$ws.Range("K37").Value = 29164.728999999999
$ws.Range("K38").Value = 525867.4
$ws.Range("K39").Value = 88922.61
$ws.Range("K40").Value = 85687.59
$ws.Range("K41").Value = 375119.9
$ws.Range("K42").Value = 32381.828000000001
$ws.Range("K43").Value = 479701.84
$ws.Range("K44").Value = 18154.312999999998
$ws.Range("K45").Value = 36187.152000000002
$ws.Range("K46").Value = 109277.71
$ws.Range("K47").Value = 90219.6
$ws.Range("K48").Value = 38230.535000000003
$ws.Range("K49").Value = 313997.03000000003
$ws.Range("K50").Value = 44513.57
$ws.Range("K51").Value = 26148.68
$ws.Range("K52").Value = 13010.584999999999
$ws.Range("K53").Value = 49924.18
$ws.Range("K54").Value = 40183.105000000003
$ws.Range("K55").Value = 2079884.8
$ws.Range("K56").Value = 20493.228999999999
$ws.Range("K57").Value = 86055.78
$ws.Range("K58").Value = 362068.38
$ws.Range("K59").Value = 5291.6934000000001
$ws.Range("K60").Value = 53592.495999999999
$ws.Range("K61").Value = 21184.440999999999
$ws.Range("K62").Value = 316975.15999999997

# Update sheet view: zoom scale and selection
$excel.ActiveWindow.Zoom = 103
$ws.Range("O5").Select()
